$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 4 de Julio de 2020 a las 09:56"

# Armenia overtook Nigeria in the ranking (rows 52/53), so the two
# country names swap places.
$ws.Range("A52").Value = "Armenia"
$ws.Range("A53").Value = "Nigeria"

# Row 6 (rank 10 - Rusia)
$ws.Range("B6").Value = 674515
$ws.Range("C6").Value = 6632
$ws.Range("D6").Value = 446879
$ws.Range("E6").Value = 217609
$ws.Range("G6").Value = 168
$ws.Range("H6").Value = 10027

# Row 7 (rank 11 - India)
$ws.Range("D7").Value = 394411
$ws.Range("E7").Value = 236809

# Row 47 (rank 51)
$ws.Range("B47").Value = 32672
$ws.Range("C47").Value = 348
$ws.Range("D47").Value = 19164
$ws.Range("E47").Value = 12682
$ws.Range("G47").Value = 7
$ws.Range("H47").Value = 826

# Row 52 (rank 56, now Armenia - updated stats that overtook Nigeria)
$ws.Range("B52").Value = 27900
$ws.Range("C52").Value = 580
$ws.Range("D52").Value = 15935
$ws.Range("E52").Value = 11488
$ws.Range("G52").Value = 8
$ws.Range("H52").Value = 477

# Row 53 (rank 57, now Nigeria - unchanged stats carried from the old row 52)
$ws.Range("B53").Value = 27564
$ws.Range("D53").Value = 11069
$ws.Range("E53").Value = 15867
$ws.Range("H53").Value = 628

# Row 96 (rank 100)
$ws.Range("B96").Value = 4174
$ws.Range("C96").Value = 2
$ws.Range("E96").Value = 834

# Row 119 (rank 123)
$ws.Range("B119").Value = 1749
$ws.Range("C119").Value = 29
$ws.Range("E119").Value = 255

# Row 132 (rank 136)
$ws.Range("B132").Value = 1123
$ws.Range("C132").Value = 1
$ws.Range("D132").Value = 1000
$ws.Range("E132").Value = 93
